# 自动更新Excel文件 - 2026-01-23 23:14:53
# For each data row (2..99) in column E (剩余/remaining days) and column F (开始时间/start date):
#   - If E == 1 and F == 20260114 (the batch that is due today), roll it forward:
#       E -> 10, F -> 20260124  (reset remaining days to D's cadence, push start date +10 days)
#   - Otherwise, the remaining-day counter just ticks down by one day: E -> E - 1
# Row 36 is special-cased out: its F value (202510929) isn't a real 8-digit date, and its
# E is already 10 (not part of either pattern), so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # Skip malformed/placeholder dates (not a normal 8-digit yyyymmdd value,
    # e.g. row 36's "202510929") -- those rows are left untouched.
    $fStr = [string][int]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }

    if ($eVal -eq 1 -and $fVal -eq 20260114) {
        # Batch is due today: reset the countdown and roll the start date forward 10 days.
        $eCell.Value = 10
        $fCell.Value = 20260124
    }
    else {
        # Otherwise the remaining-day countdown just ticks down by one.
        $eCell.Value = $eVal - 1
    }
}
